$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new data row is being inserted at sheet row 185 ("Hortaliza, Terminal
# Hortofrutícola Agro Chillán - Zapallo italiano"). Every existing row from
# 185 through 221 shifts down by one (to 186..222). We implement the shift
# manually (bottom-up, so we never clobber a row before it has been copied),
# then overwrite row 185 with the new record's values.

$firstRow = 185
$lastRow = 221
$lastCol = 18   # column R

$dateCol = 4   # column D ("Fecha") carries a date NumberFormat

for ($r = $lastRow; $r -ge $firstRow; $r--) {
    $srcRow = $r
    $dstRow = $r + 1
    for ($c = 1; $c -le $lastCol; $c++) {
        $srcCell = $ws.Cells.Item($srcRow, $c)
        $dstCell = $ws.Cells.Item($dstRow, $c)
        $dstCell.Value = $srcCell.Value2
        if ($c -eq $dateCol) {
            $dstCell.NumberFormat = $srcCell.NumberFormat
        }
    }
}

# Now populate the new row 185 with the inserted record's data.
$ws.Cells.Item(185, 1).Value = 7
$ws.Cells.Item(185, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(185, 3).Value = "Ñuble"
$ws.Cells.Item(185, 4).Value = 44694
$ws.Cells.Item(185, 4).NumberFormat = $ws.Cells.Item(186, 4).NumberFormat
$ws.Cells.Item(185, 5).Value = 16
$ws.Cells.Item(185, 6).Value = 100112032
$ws.Cells.Item(185, 7).Value = "Zapallo italiano"
$ws.Cells.Item(185, 8).Value = "Sin especificar"
$ws.Cells.Item(185, 9).Value = "Primera"
$ws.Cells.Item(185, 10).Value = 100
$ws.Cells.Item(185, 11).Value = 15000
$ws.Cells.Item(185, 12).Value = 15000
$ws.Cells.Item(185, 13).Value = 15000
$ws.Cells.Item(185, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(185, 15).Value = "Región del Maule"
$ws.Cells.Item(185, 16).Value = 300
$ws.Cells.Item(185, 17).Value = 50
$ws.Cells.Item(185, 18).Value = "Hortaliza"
